$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.415941596031189
$ws.Range("B1").Value = 0.6818171143531799
$ws.Range("C1").Value = 2.134377956390381
$ws.Range("D1").Value = 4.764405727386475
$ws.Range("E1").Value = 2.145371437072754
